$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "TOT.ORE PERIODO" column (F) for rows 2-8: each row sums its
# own Totale (D) with the matching detail-table total (F13:F19).
$ws.Range("F2").Formula = "=D2+F13"
$ws.Range("F3").Formula = "=D3+F14"
$ws.Range("F4").Formula = "=D4+F15"
$ws.Range("F5").Formula = "=D5+F16"
$ws.Range("F6").Formula = "=D6+F17"
$ws.Range("F7").Formula = "=D7+F18"
$ws.Range("F8").Formula = "=D8+F19"

# Match the cell style used elsewhere in those rows (style index "1" - the
# row's default format). Copy formatting only, from a cell that already
# carries that exact style, onto each new F cell.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("F2:F8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update the view: scroll so column E is the left-most visible column, and
# select F8 (matches the author's final on-screen state).
$ws.Range("F8").Select()
